$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the
# title (Heading1) paragraph "Play Country Farming for Free - Review".
# The new paragraph is Normal style, with a leading empty run, a bold
# run "Meta description" and a trailing normal run with the rest of
# the text - matching the pattern used elsewhere in this document.
# -----------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)

# Create a real paragraph break after the title so the new content
# lives in its own paragraph (does not merge into the next one).
$null = $titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Explore the features of Country Farming slot game with our review. Play for free and enjoy high winning potential, cartoon-style graphics, and special symbols.</w:t></w:r></w:p></w:body></w:document>'

$null = $metaRange.InsertXML($metaXml)

# -----------------------------------------------------------------
# Change 2: near the end of the document, remove the duplicate bold
# "Play Country Farming for Free - Review" paragraph entirely, and
# replace the text of the following italic paragraph with the new
# image-generation prompt (keeping its italic formatting and the
# leading empty run intact).
# -----------------------------------------------------------------

# Search from the bottom of the document (the title/meta-description
# text we just inserted up top is identical, so a whole-document Find
# could match the wrong / first occurrence - walk paragraphs in
# reverse instead to reliably find the trailing duplicate).
$dupTitle = $null
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Play Country Farming for Free - Review`r") {
        $dupTitle = $para
        break
    }
}

if ($dupTitle -ne $null) {
    $null = $dupTitle.Range.Delete()
}

$newPromptText = "Prompt: Create a feature image for Country Farming that showcases the game's fun and colorful nature while incorporating the Maya warrior character with glasses. The image should be in a cartoon style and feature the game's farm symbols such as animals and fruits. It should also include the game's logo. The Maya warrior can be depicted engaging in a fun activity or interacting with the symbols in some way, making the image lively and engaging to potential players."

# The italic "Explore the features..." paragraph is now the very last
# paragraph in the document - address it directly by index and
# replace its text via Range.Text (Find/Replace, which applies
# Word's smart-quotes autocorrect, is intentionally avoided so the
# straight apostrophes in the new text are preserved verbatim).
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastRange.Text = $newPromptText
